$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing date-format style (from B37) down onto B38:B41,
# and the existing time-format style (from C37:D37) onto C40:D40 / C41
# so the new rows reuse the workbook's existing cellXfs entries instead of
# minting new ones.
$ws.Range("B37").Copy()
$ws.Range("B38:B41").PasteSpecial(-4122)

$ws.Range("C37:D37").Copy()
$ws.Range("C40:D40").PasteSpecial(-4122)

$ws.Range("C37").Copy()
$ws.Range("C41").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row 38: 2022-10-18, sick day
$ws.Range("B38").Value = 44852
$ws.Range("C38").Value = "x"
$ws.Range("D38").Value = "x"
$ws.Range("E38").Value = "Sick"

# Row 39: 2022-10-19, sick day
$ws.Range("B39").Value = 44853
$ws.Range("C39").Value = "x"
$ws.Range("D39").Value = "x"
$ws.Range("E39").Value = "Sick"

# Row 40: 2022-10-20, back to work on pcb
$ws.Range("B40").Value = 44854
$ws.Range("C40").Value = 0.41666666666666669
$ws.Range("D40").Value = 0.66666666666666663
$ws.Range("E40").Value = "worked on pcb, talked with johan b about pcb, talked with johan k on software design"

# Row 41: 2022-10-21, started
$ws.Range("B41").Value = 44855
$ws.Range("C41").Value = 0.49305555555555558

# New left-aligned style for the "Bezigheden" notes column on these rows
$ws.Range("E38:E40").HorizontalAlignment = -4131

# Update view to match author's scroll/selection position
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("E40").Select()
